$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.709.27"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.351.07"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -11.70%  "
$ws.Range("D9").Value = "2.349.60"
$ws.Range("E9").Value = "  -4.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.95%  "
$ws.Range("D15").Value = "2.771.79"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").Value = "60.417.44"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.348.49"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "2.462.87"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").Value = "0.0₃0890"
$ws.Range("E30").Value = "  -9.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "501.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.17%  "
$ws.Range("E32").Value = "  -5.64%  "
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  -6.00%  "
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.371"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.12%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.569"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("E51").Value = "  -4.49%  "
